$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.7511267625674056
$ws.Range("C2").Value = 1.174505337654612
$ws.Range("D2").Value = 3.069342633011421
$ws.Range("E2").Value = 0.7851754730668936
$ws.Range("D3").Value = 2.288657025972322
$ws.Range("E3").Value = 0.8392668032967509
$ws.Range("B4").Value = 0.6947032220235132
$ws.Range("C4").Value = 1.510862277990596
$ws.Range("D4").Value = -0.1663351976388032
$ws.Range("E4").Value = 0.953330387505735
$ws.Range("D5").Value = 1.743425413331225
$ws.Range("E5").Value = 0.4709357672430788
$ws.Range("B6").Value = 0.54790035230731
$ws.Range("E6").Value = 0.5479003523073096
$ws.Range("B7").Value = 0.4350736390815085
$ws.Range("C7").Value = 0.2023621702127827
$ws.Range("D7").Value = 0.3979485974093362
$ws.Range("E7").Value = 0.8793012627533111
$ws.Range("B8").Value = 0.5982023189954352
$ws.Range("C8").Value = 1.889282499519986
$ws.Range("D8").Value = 5.246821010190586
$ws.Range("E8").Value = 0.6632706282838209
$ws.Range("D9").Value = 4.085552086025032
$ws.Range("E9").Value = 0.7381685387562538
$ws.Range("B10").Value = 0.6269980998833792
$ws.Range("C10").Value = 1.661284310198333
$ws.Range("D10").Value = 2.117509708730804
$ws.Range("E10").Value = 0.798964843940331
$ws.Range("B11").Value = 0.3718699914299579
$ws.Range("C11").Value = 75.98583345497204
$ws.Range("D11").Value = 261.7080780761954
$ws.Range("E11").Value = 0.4744101715212642
$ws.Range("B12").Value = 0.5476973011871676
$ws.Range("D12").Value = 227.5188879085707
$ws.Range("E12").Value = 0.5476973011871671
$ws.Range("B13").Value = 0.4515593679468755
$ws.Range("C13").Value = 69.83501878436971
$ws.Range("D13").Value = 140.3680615770823
$ws.Range("E13").Value = 0.765628649361459
$ws.Range("B14").Value = 0.4357359037190183
$ws.Range("D14").Value = 0.4631455524918361
$ws.Range("B15").Value = 0.6991188334340351
$ws.Range("C15").Value = 0.2536382222984022
$ws.Range("D15").Value = 0.3117524329442313
$ws.Range("E15").Value = 0.699118833434035
$ws.Range("B16").Value = 0.3055639137124027
$ws.Range("C16").Value = 0.5567967778383045
$ws.Range("D16").Value = 0.08942439561234861
$ws.Range("E16").Value = 1.168541207869414
